$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad / "Changed date") for rows 2-10 changes from
# serial date 45208 (2023-10-09) to 45212 (2023-10-13).
$ws.Range("C2:C10").Value = 45212
